$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 2021 row (row 2) entirely; rows 3-6 shift up to rows 2-5.
$ws.Rows.Item(2).Delete()

# Update retention_rate for the (now first) 2022 data row, which was
# recalculated now that there is no preceding year to compare against.
$ws.Range("F2").Value = 35.0104821802935
